$d = $word.ActiveDocument

$table = $d.Tables.Item(1)

$map = @{
    "1,1" = "20÷3="
    "1,2" = "25÷2="
    "1,3" = "95÷9="
    "1,4" = "93÷7="
    "1,5" = "80÷6="
    "5,1" = "12÷8="
    "5,2" = "39÷9="
    "5,3" = "22÷9="
    "5,4" = "61÷4="
    "5,5" = "16÷4="
    "9,1" = "37÷4="
    "9,2" = "97÷3="
    "9,3" = "43÷3="
    "9,4" = "64÷2="
    "9,5" = "90÷4="
    "13,1" = "78÷7="
    "13,2" = "22÷5="
    "13,3" = "93÷3="
    "13,4" = "52÷3="
    "13,5" = "48÷5="
    "17,1" = "97÷9="
    "17,2" = "83÷8="
    "17,3" = "79÷7="
    "17,4" = "12÷5="
    "17,5" = "87÷8="
}

foreach ($key in $map.Keys) {
    $parts = $key.Split(",")
    $row = [int]$parts[0]
    $col = [int]$parts[1]
    $cell = $table.Cell($row, $col)
    $cellRange = $cell.Range
    $cellRange.End = $cellRange.End - 1
    $cellRange.Text = $map[$key]
}
